$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2 = '  -2.06%  '
    3 = '  -4.91%  '
    4 = '  +0.01%  '
    5 = '  -1.52%  '
    6 = '  -7.56%  '
    7 = '  +0.04%  '
    8 = '  -10.86%  '
    9 = '  -4.83%  '
    10 = '  -2.15%  '
    11 = '  -0.03%  '
    12 = '  -2.38%  '
    13 = '  -3.49%  '
    14 = '  -6.54%  '
    15 = '  -4.91%  '
    16 = '  -1.64%  '
    17 = '  -4.72%  '
    18 = '  -4.89%  '
    19 = '  -4.12%  '
    20 = '  -1.08%  '
    21 = '  -3.02%  '
    22 = '  -6.11%  '
    23 = '  -0.20%  '
    24 = '  -2.99%  '
    25 = '  -1.36%  '
    26 = '  +11.17%  '
    27 = '  +0.22%  '
    28 = '  -4.84%  '
    29 = '  -9.93%  '
    30 = '  -4.16%  '
    31 = '  -6.21%  '
    32 = '  -7.00%  '
    33 = '  -1.73%  '
    34 = '  -5.07%  '
    35 = '  -5.29%  '
    36 = '  +0.06%  '
    37 = '  -4.31%  '
    38 = '  -2.13%  '
    39 = '  -0.19%  '
    40 = '  -9.29%  '
    41 = '  +1.95%  '
    42 = '  -0.06%  '
    43 = '  -2.36%  '
    44 = '  -0.34%  '
    45 = '  -1.41%  '
    46 = '  -8.23%  '
    47 = '  -2.53%  '
    48 = '  -4.62%  '
    49 = '  -12.16%  '
    50 = '  -4.09%  '
    51 = '  -3.89%  '
}

foreach ($row in $values.Keys) {
    $ws.Range("E$row").Value = $values[$row]
}
